# Applies the "Changed values in testData" commit to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: password/value updates ---
$ws.Range("B1").Value = "Payable"
$ws.Range("C1").Value = "Testing1@"

# --- Row 4: vendor name + account updates ---
$ws.Range("B4").Value = "TechBite"
$ws.Range("Q4").Value = "Account"

# --- Row 5: new invoice test data row (was previously just A5 populated) ---
$ws.Range("B5").Value = "TechBite"
$ws.Range("C5").Value = "Net 30"
$ws.Range("D5").Value = "pune"
$ws.Range("E5").Value = "Advertising"
$ws.Range("F5").Value = "invoice desc 1"
$ws.Range("G5").Value = 500
$ws.Range("H5").Value = "Workbooks"
$ws.Range("I5").Value = "Department1"
$ws.Range("J5").Value = "Advertising"
$ws.Range("K5").Value = "inv desc2"
$ws.Range("L5").Value = "measure2"
$ws.Range("M5").Value = 234.44444999999999
$ws.Range("N5").Value = 876.9757366

# --- Row 6: new invoice test data row (was previously just A6 populated) ---
$ws.Range("B6").Value = "TechBite"
$ws.Range("C6").Value = "Net 30"
$ws.Range("D6").Value = "pune"
$ws.Range("E6").Value = "Advertising"
$ws.Range("F6").Value = "invoice desc 1"
$ws.Range("G6").Value = 500
$ws.Range("H6").Value = "laptop"
$ws.Range("I6").Value = "Department1"
$ws.Range("J6").Value = "Bank Charges"
$ws.Range("K6").Value = "inv desc3"
$ws.Range("L6").Value = "measure3"
$ws.Range("M6").Value = 38.741129999999998
$ws.Range("N6").Value = 34.987736650000002

# --- Row 9: AP Vendor header label update ---
$ws.Range("AC9").Value = "NetChain2 – AP Vendor: New"

# --- Row 13: new sales order test data row ---
$ws.Range("A13").Value = "AR.NetchainTest.CreateSalesOrder"
$ws.Range("B13").Value = "Accenture"
$ws.Range("C13").Value = "pune"
$ws.Rows.Item(13).RowHeight = 35.25

# --- View state: widen column AC (stored width 27) and move the selection/scroll position ---
$ws.Columns.Item(29).ColumnWidth = 26.14
$ws.Range("Q4").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 13
